# Apply the "Add files via upload" revision to the soil-profile worksheet.
#
# Content changes (SampleNo / Soil Class columns):
#   - A15: "DB-12" -> "SS-12"
#   - A16: "DB-13" -> "SS-12"
#   - K15: "SC"    -> "Granite"
#   - K16: "SC"    -> "Decomposed Granite"
#   - K17: "SC"    -> "SC/Decomposed Granite"
#   - K18: "SC"    -> "SC/Decomposed Granite"
#   - K19: "SC"    -> "SC/Decomposed Granite"
#   - K20: "SC"    -> "(SM/Decomposed Granite)"
#   - K21: "SC"    -> "(SM/Decomposed Granite)"
#   - K22: "SC"    -> "Granite"
#   - K23: "SC"    -> "Granite"
#   - K24: "SC"    -> "Granite"
#   - K25: "SC"    -> "Granite"
#
# View changes: the sheet is now scrolled/zoomed differently and a new
# cell is selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates -------------------------------------------------

$ws.Range("A15").Value = "SS-12"
$ws.Range("A16").Value = "SS-12"

$ws.Range("K15").Value = "Granite"
$ws.Range("K16").Value = "Decomposed Granite"
$ws.Range("K17").Value = "SC/Decomposed Granite"
$ws.Range("K18").Value = "SC/Decomposed Granite"
$ws.Range("K19").Value = "SC/Decomposed Granite"
$ws.Range("K20").Value = "(SM/Decomposed Granite)"
$ws.Range("K21").Value = "(SM/Decomposed Granite)"
$ws.Range("K22").Value = "Granite"
$ws.Range("K23").Value = "Granite"
$ws.Range("K24").Value = "Granite"
$ws.Range("K25").Value = "Granite"

# --- View / selection updates ---------------------------------------------
# topLeftCell moves from E1 to A3, zoom goes from 142% to 259%, and the
# active selection moves from L20 to A15.

$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$win.Zoom = 259

$ws.Range("A15").Select()
